$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete column-A labels before re-using their rows.
#     Deleted highest row first so the up-shift never disturbs a
#     not-yet-processed cell. Neither A26 nor A31 has anything
#     below it in column A, so this is a clean removal. ---
$ws.Range("A31").Delete()
$ws.Range("A26").Delete()

# --- Row 17: Float Şarj section header (text of B17 changes) ---
$ws.Range("A17").Value = "Float Şarj"
$ws.Range("B17").Value = "Ayarlardan bu moda ait voltaj ve akım değerini oku ve yaz."

# --- New row 18: Boost Şarj header + "Bu modu devreye al." ---
$ws.Range("A18").Value = "Boost Şarj"
$ws.Range("B18").Value = "Bu modu devreye al."

# --- New row 19: "Testleri" + "Batarya bağla. Hafif yük ver." ---
$ws.Range("A19").Value = "Testleri"
$ws.Range("B19").Value = "Batarya bağla. Hafif yük ver."

# --- Rows 20-28: updated test-step instructions ---
$ws.Range("B20").Value = "Şarj voltajını kontrol et. Ön panelden oku, ölçüm cihazı ile ölç ve yaz."
$ws.Range("B21").Value = "Yükü artırarak akım sınırına gel."
$ws.Range("B22").Value = "Akım sınırını kontrol et. Ön panelden oku, ölçüm cihazı ile ölç ve yaz."
$ws.Range("B23").Value = "Akım ve voltajı değiştirme menüsüne gel."
$ws.Range("B24").Value = "Akım ve voltajı değiştir ve yaz."
$ws.Range("B25").Value = "Tekrar akım ve voltajı kontrol et ve yaz."
$ws.Range("B26").Value = "Ölçüm değeri ile ayar değerini karşılaştır."
$ws.Range("B27").Value = "Cihazı kapatıp aç. Ayar değeri duruyor mu kontrol et."
$ws.Range("B28").Value = "Akım ve voltajı ölç ve yaz."

# --- New rows 29-30: Otomatik Şarj section ---
$ws.Range("A29").Value = "Otomatik Şarj"
$ws.Range("B29").Value = "Float ve boost testleri tamamlanmış varsayılıyor."
$ws.Range("B30").Value = "Manuel olarak float şarj moduna geç."

# --- "Oto şarj test" moves from A31 down to A33 ---
$ws.Range("A33").Value = "Oto şarj test"

# --- Apply background fill to the Float/Boost test-step block ---
# (only the cells that actually hold content get the fill: A17:B19 plus B20:B28)
$ws.Range("A17:B19").Interior.Color = 15722206
$ws.Range("B20:B28").Interior.Color = 15722206

# --- Widen column B to fit the longer instructions ---
$ws.Columns("B").ColumnWidth = 62

# --- Restore a sensible selection/scroll position ---
$ws.Range("B31").Select()
